# measure_results all but the wins calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Won" column (D) is set equal to the "Games" column (C) for data rows 2-9
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $ws.Cells.Item($r, 3).Value2
}

# "Percent" column (E) becomes "100%" for every data row.
# Row 8 (E8) already holds the text "100%", so copy that value (not the
# format) onto the other Percent cells to avoid introducing any new
# number format / style.
$ws.Range("E8").Copy()
for ($r = 2; $r -le 9; $r++) {
    if ($r -ne 8) {
        $ws.Range("E$r").PasteSpecial(-4163) # xlPasteValues
    }
}
$excel.CutCopyMode = $false
